$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before H ("Run" status column) - this shifts the
# existing H:M (task_id..NumberOfMissingValues) columns right to I:N.
$ws.Columns("H").Insert()

# Header for the newly inserted column.
$ws.Range("H1").Value = "MA2"
$ws.Columns("H").ColumnWidth = 16.83

# Populate the new "Run" status column: wherever the row already has a
# value in the G ("Status") column, mark the new H column as "Run" too.
for ($r = 2; $r -le 93; $r++) {
    $g = $ws.Cells.Item($r, 7).Value()
    if ($g -ne $null) {
        $ws.Cells.Item($r, 8).Value = "Run"
    }
}

# Restore the active selection on the newly added column's header cell
# (the frozen-pane top-left cell is reset to B2 automatically by the
# column insert above, matching the target view state).
$null = $ws.Range("H1").Select()
